# Generate Report for Handoff
#
# The "b.md" file just went through a fresh handoff cycle (new xlf files
# were generated for it), so:
#   - Overview sheet: its per-language status cells move from
#     "Handed back: in sync with en-US" to "Ready for handoff", with an
#     updated "Latest HO Xliff Generate Date".
#   - zh-cn / de-de detail sheets: the a.md row's Status text is refreshed
#     to the same new wording, and the b.md row gets the new handoff file
#     name/date, its "Content Duplicate" flag flips to False (as text, to
#     match the existing column formatting) because the content changed,
#     and an Error Detail note is added flagging that the previously
#     recorded handback file is now out of date.

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"
$overviewDate    = "2016-09-06 04:47:42"
$staleWarning    = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b74bb9556b262c1e4889ccfb979c5f417d45a65e/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/21a5190cabeb6a5fcc31fd83f46b8a122f7c4bc4/e2e/b.md."

# Helper: write a literal text value even when it looks like a boolean
# keyword ("True"/"False") so Excel's auto-type-detection doesn't silently
# turn it into a real boolean cell. The leading "'" forces text entry (and
# is stripped from the stored content); re-applying the Normal style then
# clears the resulting quote-prefix flag so the cell's formatting matches
# a plain, never-touched text cell.
function Set-TextValue($range, $text) {
    $range.Formula = "'" + $text
    $range.Style = "Normal"
}

# --- Overview sheet: b.md row (row 3) -------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $readyForHandoff
$wsOverview.Range("F3").Value = $readyForHandoff
$wsOverview.Range("G3").Value = $overviewDate

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# a.md row: status wording refreshed too
$wsZhCn.Range("C2").Value = $readyForHandoff

# b.md row: new handoff generated, old handback now stale
$wsZhCn.Range("C3").Value = $readyForHandoff
Set-TextValue $wsZhCn.Range("F3") "False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-06 04:47:37"
$wsZhCn.Range("P3").Value = $staleWarning

# Error Detail column now holds long text -> widen it like the other
# file-name columns (G / J), which are already 40 wide.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet --------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $readyForHandoff

$wsDeDe.Range("C3").Value = $readyForHandoff
Set-TextValue $wsDeDe.Range("F3") "False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = $overviewDate
$wsDeDe.Range("P3").Value = $staleWarning

$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
